$d = $word.ActiveDocument

# 1. Refresh the footer timestamp.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("2025-06-30 12:13Z / ", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# 2. Regression test: add the standard HTML-ish character styles (b, i, sub, sup, u).
$styleDefs = @(
    @{ Id = "b";   Bold = $true },
    @{ Id = "i";   Italic = $true },
    @{ Id = "sub"; VerticalAlign = "subscript" },
    @{ Id = "sup"; VerticalAlign = "superscript" },
    @{ Id = "u";   Underline = "single" }
)

foreach ($def in $styleDefs) {
    $s = $d.Styles.Add($def.Id, 2)
    $s.BaseStyle = "DefaultParagraphFont"
    $s.Priority = 1
    $s.QuickStyle = $true

    if ($def.ContainsKey("Bold")) {
        $s.Font.Bold = $def.Bold
    }
    if ($def.ContainsKey("Italic")) {
        $s.Font.Italic = $def.Italic
    }
    if ($def.ContainsKey("VerticalAlign")) {
        if ($def.VerticalAlign -eq "subscript") {
            $s.Font.Subscript = $true
        } elseif ($def.VerticalAlign -eq "superscript") {
            $s.Font.Superscript = $true
        }
    }
    if ($def.ContainsKey("Underline")) {
        $s.Font.Underline = 1
    }
}
